# Swap the order of slides 9 and 10 (both titled "Declaração de variáveis e
# constantes"): the slide that currently sits at position 10 (the shorter,
# "Antes de escrevermos..." version) moves to position 9, and the slide
# that currently sits at position 9 (the longer, "Para definir uma
# variável..." version with the brace/key shapes) shifts down to position 10.
$p = $ppt.ActivePresentation
$p.Slides.Item(10).MoveTo(9)
